# DC-Colos.xlsx — add a new colocation-facility row for Xining, China (XNN)
# right before the existing IAD (Ashburn, VA) row. Every row from the old
# row 271 (IAD) through the old row 330 (YHZ / Halifax) shifts down by one.
#
# The new row has no lat/long data (same situation as several other
# newly-added Asia rows in this sheet, e.g. TEN/Tongren and HYN/Taizhou
# immediately above it), so columns G and H are left as empty text cells
# rather than numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push IAD..YHZ (old rows 271-330) down one row, opening up row 271.
$ws.Rows(271).Insert()

# Row 270 (HYN / Taizhou) already has the right look for a freshly-added
# Asia row (bold/centered/bordered colo code in column A, plain text for
# the rest, and blank G/H) — copy its formatting onto the new row 271
# before filling in the real values.
$ws.Range("A270:H270").Copy()
$ws.Range("A271:H271").PasteSpecial(-4122)

$ws.Range("A271").Value = "XNN"
$ws.Range("B271").Value = "Xining, China"
$ws.Range("C271").Value = "Asia"
$ws.Range("D271").Value = "Xining"
$ws.Range("E271").Value = "China"
$ws.Range("F271").Value = "CN"

# G271/H271: no coordinates for this facility. Write a bare quote-prefix
# (forces an explicit, empty text value instead of leaving the cell a
# untyped/blank) then strip the auto-applied "quoted text" style back to
# Normal so the cells match the plain, unstyled empty-text cells used
# elsewhere in the sheet for colos without lat/long data.
$ws.Range("G271").Value = "'"
$ws.Range("H271").Value = "'"
$ws.Range("G271").Style = "Normal"
$ws.Range("H271").Style = "Normal"
